# Add SpoVM, Magainin2, GMAP210 to AH list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended after the last existing data row (row 42 -> 43,44,45):
#   row43: SpoVM     / MKFYTIKLPKFLGGIVRAMLGSFRKD
#   row44: Magainin2 / GIGKFLHSAKKFGKAFVGEIMNS
#   row45: GMAP210   / MSSWLGGLGSGLGQSLGQVGGSLASLTGQISNFTKDML
# Values are entered in the same order the author typed them so the
# resulting shared-string table ordering matches the source workbook.
$ws.Range("B44").Value = "Magainin2"
$ws.Range("B43").Value = "SpoVM"
$ws.Range("B45").Value = "GMAP210"

$ws.Range("A43").Value = "MKFYTIKLPKFLGGIVRAMLGSFRKD"
$ws.Range("A45").Value = "MSSWLGGLGSGLGQSLGQVGGSLASLTGQISNFTKDML"
$ws.Range("A44").Value = "GIGKFLHSAKKFGKAFVGEIMNS"

$ws.Range("C43").Value = 1
$ws.Range("C44").Value = 1
$ws.Range("C45").Value = 1

# Match final view state from the diff: scrolled so row 23 is at top, active cell A44
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("A44").Select()
